$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-10 23:48:37'
$ws.Range('O2').Value = '1.0 °C'
$ws.Range('E3').Value = '2026-02-10 23:48:40'
$ws.Range('O3').Value = '1.0 °C'
$ws.Range('E4').Value = '2026-02-10 23:48:42'
$ws.Range('E5').Value = '2026-02-10 23:48:45'
$ws.Range('E6').Value = '2026-02-10 23:48:47'
$ws.Range('J6').Value = '1003.9 hPa'
$ws.Range('O6').Value = '9.9 °C'
$ws.Range('E7').Value = '2026-02-10 23:48:49'
$ws.Range('H7').Value = "'69%"
$ws.Range('J7').Value = '1004.1 hPa'
$ws.Range('E8').Value = '2026-02-10 23:48:52'
$ws.Range('J8').Value = '1004.0 hPa'
$ws.Range('E9').Value = '2026-02-10 23:48:55'
$ws.Range('E10').Value = '2026-02-10 23:48:57'
$ws.Range('E11').Value = '2026-02-10 23:48:59'
$ws.Range('H11').Value = "'91%"
$ws.Range('E12').Value = '2026-02-10 23:49:02'
$ws.Range('E13').Value = '2026-02-10 23:49:04'
$ws.Range('H13').Value = "'92%"
$ws.Range('J13').Value = '1006.3 hPa'
$ws.Range('E14').Value = '2026-02-10 23:49:07'
$ws.Range('O14').Value = '12.9 °C'
$ws.Range('E15').Value = '2026-02-10 23:49:09'
$ws.Range('E16').Value = '2026-02-10 23:49:11'
$ws.Range('H16').Value = "'85%"
$ws.Range('E17').Value = '2026-02-10 23:49:14'
$ws.Range('L17').Value = '135.0 km/h - 267º 23:25 TU'
$ws.Range('E18').Value = '2026-02-10 23:49:16'
$ws.Range('E19').Value = '2026-02-10 23:49:19'
$ws.Range('E20').Value = '2026-02-10 23:49:21'
$ws.Range('I20').Value = '12.3 mm'
$ws.Range('O20').Value = '0.7 °C'
$ws.Range('E21').Value = '2026-02-10 23:49:24'
$ws.Range('E22').Value = '2026-02-10 23:49:27'
$ws.Range('O22').Value = '-0.5 °C'
$ws.Range('E23').Value = '2026-02-10 23:49:29'
$ws.Range('I23').Value = '27.2 mm'
$ws.Range('E24').Value = '2026-02-10 23:49:32'
$ws.Range('L24').Value = '25.9 km/h - 277º 23:29 TU'
$ws.Range('E25').Value = '2026-02-10 23:49:34'
$ws.Range('E26').Value = '2026-02-10 23:49:36'
$ws.Range('E27').Value = '2026-02-10 23:49:38'
$ws.Range('L27').Value = '58.0 km/h - 227º 23:05 TU'
$ws.Range('O27').Value = '1.4 °C'
$ws.Range('E28').Value = '2026-02-10 23:49:41'
$ws.Range('E29').Value = '2026-02-10 23:49:44'
$ws.Range('E30').Value = '2026-02-10 23:49:46'
$ws.Range('E31').Value = '2026-02-10 23:49:48'
$ws.Range('H31').Value = "'81%"
$ws.Range('L31').Value = '52.6 km/h - 326º 23:26 TU'
$ws.Range('M31').Value = '15.7 °C 23:25 TU'
$ws.Range('O31').Value = '10.8 °C'
$ws.Range('E32').Value = '2026-02-10 23:49:51'
$ws.Range('E33').Value = '2026-02-10 23:49:54'
$ws.Range('J33').Value = '1006.0 hPa'
$ws.Range('O33').Value = '4.5 °C'
$ws.Range('E34').Value = '2026-02-10 23:49:57'
$ws.Range('E35').Value = '2026-02-10 23:49:59'
$ws.Range('O35').Value = '13.0 °C'
$ws.Range('E36').Value = '2026-02-10 23:50:01'
$ws.Range('E37').Value = '2026-02-10 23:50:04'
$ws.Range('E38').Value = '2026-02-10 23:50:06'
$ws.Range('E39').Value = '2026-02-10 23:50:09'
$ws.Range('E40').Value = '2026-02-10 23:50:12'
$ws.Range('H40').Value = "'92%"
$ws.Range('O40').Value = '7.5 °C'
$ws.Range('E41').Value = '2026-02-10 23:50:14'
$ws.Range('H41').Value = "'77%"
$ws.Range('J41').Value = '1004.2 hPa'
$ws.Range('L41').Value = '69.5 km/h - 294º 23:10 TU'
$ws.Range('O41').Value = '15.3 °C'
$ws.Range('E42').Value = '2026-02-10 23:50:17'
$ws.Range('E43').Value = '2026-02-10 23:50:20'
$ws.Range('O43').Value = '10.3 °C'
$ws.Range('E44').Value = '2026-02-10 23:50:22'
$ws.Range('E45').Value = '2026-02-10 23:50:25'
$ws.Range('J45').Value = '1005.5 hPa'
$ws.Range('E46').Value = '2026-02-10 23:50:28'
$ws.Range('O46').Value = '15.4 °C'
